# Insert a new data row at row 326 (pushing the existing rows 326:398 down to
# 327:399) and populate the new row with a fresh "Poroto verde" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 326:398 down by one row.
$ws.Rows.Item(326).Insert()

# Fill in the newly inserted row 326 with the new record's values.
$ws.Cells.Item(326, 1).Value = 3
$ws.Cells.Item(326, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(326, 3).Value = "Coquimbo"
$ws.Cells.Item(326, 4).Value = 44754
$ws.Cells.Item(326, 5).Value = 5
$ws.Cells.Item(326, 6).Value = 100112031
$ws.Cells.Item(326, 7).Value = "Poroto verde"
$ws.Cells.Item(326, 8).Value = "Magnum"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 82
$ws.Cells.Item(326, 11).Value = 33000
$ws.Cells.Item(326, 12).Value = 34000
$ws.Cells.Item(326, 13).Value = 33488
$ws.Cells.Item(326, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(326, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(326, 16).Value = 1340
$ws.Cells.Item(326, 17).Value = 25
$ws.Cells.Item(326, 18).Value = "Hortaliza"

# Give the new date cell (D326) the same date style as the rest of column D.
$ws.Cells.Item(326, 4).NumberFormat = $ws.Cells.Item(327, 4).NumberFormat
